# BBI-23-6.xlsx: the pass/fail threshold on column K ("прошел отбор") was
# loosened from 57.5/4 to 57.5/5, so more rows now flag as passing (1).
#
# K2 holds its own literal formula; K3:K26 is one shared formula group
# (master cell K3). Re-write both so the shared group keeps working the
# same way it did before (Excel will re-infer the shared range itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("K2").Formula = "=IF(J2>57.5/5,1,0)"
$ws.Range("K3:K26").Formula = "=IF(J3>57.5/5,1,0)"

# Leave the cursor where the author left it after editing the formula:
# active cell K2, with the whole dependent column selected.
$ws.Range("K2:K26").Select()
